$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '67.701.96'
$ws.Range('E2').Value = '  -1.99%  '

# Row 3
$ws.Range('D3').Value = '2.679.22'
$ws.Range('E3').Value = '  -2.26%  '

# Row 4
$ws.Range('E4').Value = '  +0.03%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '600.91'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.40%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '167.53'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.50%  '

# Row 7
$ws.Range('E7').Value = '  +0.02%  '

# Row 8
$ws.Range('E8').Value = '  +0.07%  '

# Row 9
$ws.Range('D9').Value = '2.679.11'
$ws.Range('E9').Value = '  -2.27%  '

# Row 10
$ws.Range('E10').Value = '  +1.30%  '

# Row 11
$ws.Range('E11').Value = '  +1.11%  '

# Row 12
$ws.Range('E12').Value = '  +0.04%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.23'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.79%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.93'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.76%  '

# Row 15
$ws.Range('D15').Value = '3.165.02'
$ws.Range('E15').Value = '  -2.38%  '

# Row 16
$ws.Range('E16').Value = '  -2.85%  '

# Row 17
$ws.Range('D17').Value = '67.953.15'
$ws.Range('E17').Value = '  -1.45%  '

# Row 18
$ws.Range('D18').Value = '2.680.86'
$ws.Range('E18').Value = '  -1.59%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.75'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.70%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.90'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.10%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '365.29'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.89%  '

# Row 22
$ws.Range('E22').Value = '  -3.58%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.83'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.88%  '

# Row 24
$ws.Range('E24').Value = '  -4.93%  '

# Row 25
$ws.Range('E25').Value = '  +0.07%  '

# Row 26
$ws.Range('E26').Value = '  -4.21%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.27'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.33%  '

# Row 28
$ws.Range('D28').Value = '2.834.39'

# Row 29
$ws.Range('E29').Value = '  -3.78%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.999'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.01%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '554.30'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -7.97%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.01'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.92%  '

# Row 33
$ws.Range('E33').Value = '  -3.89%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.93'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.59%  '

# Row 35
$ws.Range('E35').Value = '  -1.34%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.00'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.00%  '

# Row 37
$ws.Range('E37').Value = '  -5.23%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '19.52'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -3.07%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '155.57'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -4.61%  '

# Row 40
$ws.Range('E40').Value = '  -2.56%  '

# Row 41
$ws.Range('B41').Value = 'Stacks'
$ws.Range('C41').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.84'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -4.45%  '

# Row 42
$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.31'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.49%  '

# Row 43
$ws.Range('E43').Value = '  -0.55%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.52'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -7.03%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '40.43'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.88%  '

# Row 47
$ws.Range('D47').Value = '0.0₆0301'
$ws.Range('E47').Value = '  -5.29%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.591'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.41%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '153.87'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.45%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.88'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.36%  '

# Row 51
$ws.Range('E51').Value = '  -3.93%  '
